# Update NATMI TPM metrics for Ccl11-Ccr5 sheet (rows 2-17, columns E-T)
# per refreshed TPM recompute (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = [double]"0.738254"
$ws.Range("H2").Value = [double]"2.214762"
$ws.Range("I2").Value = [double]"0.005691320045803731"
$ws.Range("J2").Value = [double]"0.005691320045803731"
$ws.Range("M2").Value = [double]"0.007957000000000001"
$ws.Range("N2").Value = [double]"0.023871"
$ws.Range("O2").Value = [double]"0.0002448939493579708"
$ws.Range("P2").Value = [double]"0.0002448939493579708"
$ws.Range("Q2").Value = [double]"0.005874287078"
$ws.Range("R2").Value = [double]"0.052868583702"
$ws.Range("S2").Value = [double]"1.393769843077063E-06"
$ws.Range("T2").Value = [double]"1.393769843077063E-06"

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = [double]"0.738254"
$ws.Range("H3").Value = [double]"2.214762"
$ws.Range("I3").Value = [double]"0.005691320045803731"
$ws.Range("J3").Value = [double]"0.005691320045803731"
$ws.Range("O3").Value = [double]"0.003249135679578298"
$ws.Range("P3").Value = [double]"0.003249135679578299"
$ws.Range("Q3").Value = [double]"0.07793722869533333"
$ws.Range("R3").Value = [double]"0.701435058258"
$ws.Range("S3").Value = [double]"1.84918710247201E-05"
$ws.Range("T3").Value = [double]"1.84918710247201E-05"

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = [double]"0.738254"
$ws.Range("H4").Value = [double]"2.214762"
$ws.Range("I4").Value = [double]"0.005691320045803731"
$ws.Range("J4").Value = [double]"0.005691320045803731"
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = [double]"0.3333333333333333"
$ws.Range("M4").Value = [double]"0.02449766666666667"
$ws.Range("N4").Value = [double]"0.073493"
$ws.Range("O4").Value = [double]"0.0007539688752111494"
$ws.Range("P4").Value = [double]"0.0007539688752111494"
$ws.Range("Q4").Value = [double]"0.01808550040733333"
$ws.Range("R4").Value = [double]"0.162769503666"
$ws.Range("S4").Value = [double]"4.291078173401306E-06"
$ws.Range("T4").Value = [double]"4.291078173401306E-06"

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = [double]"0.738254"
$ws.Range("H5").Value = [double]"2.214762"
$ws.Range("I5").Value = [double]"0.005691320045803731"
$ws.Range("J5").Value = [double]"0.005691320045803731"
$ws.Range("M5").Value = [double]"32.353591"
$ws.Range("N5").Value = [double]"97.060773"
$ws.Range("O5").Value = [double]"0.9957520014958525"
$ws.Range("P5").Value = [double]"0.9957520014958525"
$ws.Range("Q5").Value = [double]"23.885167970114"
$ws.Range("R5").Value = [double]"214.966511731026"
$ws.Range("S5").Value = [double]"0.005667143326762532"
$ws.Range("T5").Value = [double]"0.005667143326762532"

# Row 6
$ws.Range("I6").Value = [double]"0.9440493064670392"
$ws.Range("J6").Value = [double]"0.9440493064670391"
$ws.Range("M6").Value = [double]"0.007957000000000001"
$ws.Range("N6").Value = [double]"0.023871"
$ws.Range("O6").Value = [double]"0.0002448939493579708"
$ws.Range("P6").Value = [double]"0.0002448939493579708"
$ws.Range("Q6").Value = [double]"0.9743990141730002"
$ws.Range("R6").Value = [double]"8.769591127557002"
$ws.Range("S6").Value = [double]"0.0002311919630493665"
$ws.Range("T6").Value = [double]"0.0002311919630493665"

# Row 7
$ws.Range("I7").Value = [double]"0.9440493064670392"
$ws.Range("J7").Value = [double]"0.9440493064670391"
$ws.Range("O7").Value = [double]"0.003249135679578298"
$ws.Range("P7").Value = [double]"0.003249135679578299"
$ws.Range("S7").Value = [double]"0.003067344284923205"
$ws.Range("T7").Value = [double]"0.003067344284923205"

# Row 8
$ws.Range("I8").Value = [double]"0.9440493064670392"
$ws.Range("J8").Value = [double]"0.9440493064670391"
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.02449766666666667"
$ws.Range("N8").Value = [double]"0.073493"
$ws.Range("O8").Value = [double]"0.0007539688752111494"
$ws.Range("P8").Value = [double]"0.0007539688752111494"
$ws.Range("Q8").Value = [double]"2.999937444959"
$ws.Range("R8").Value = [double]"26.999437004631"
$ws.Range("S8").Value = [double]"0.0007117837937408192"
$ws.Range("T8").Value = [double]"0.0007117837937408191"

# Row 9
$ws.Range("I9").Value = [double]"0.9440493064670392"
$ws.Range("J9").Value = [double]"0.9440493064670391"
$ws.Range("M9").Value = [double]"32.353591"
$ws.Range("N9").Value = [double]"97.060773"
$ws.Range("O9").Value = [double]"0.9957520014958525"
$ws.Range("P9").Value = [double]"0.9957520014958525"
$ws.Range("Q9").Value = [double]"3961.9589261476"
$ws.Range("R9").Value = [double]"35657.63033532839"
$ws.Range("S9").Value = [double]"0.9400389864253258"
$ws.Range("T9").Value = [double]"0.9400389864253257"

# Row 10
$ws.Range("G10").Value = [double]"5.698467"
$ws.Range("H10").Value = [double]"17.095401"
$ws.Range("I10").Value = [double]"0.0439304080539368"
$ws.Range("J10").Value = [double]"0.04393040805393679"
$ws.Range("M10").Value = [double]"0.007957000000000001"
$ws.Range("N10").Value = [double]"0.023871"
$ws.Range("O10").Value = [double]"0.0002448939493579708"
$ws.Range("P10").Value = [double]"0.0002448939493579708"
$ws.Range("Q10").Value = [double]"0.045342701919"
$ws.Range("R10").Value = [double]"0.408084317271"
$ws.Range("S10").Value = [double]"1.075829112523579E-05"
$ws.Range("T10").Value = [double]"1.075829112523579E-05"

# Row 11
$ws.Range("G11").Value = [double]"5.698467"
$ws.Range("H11").Value = [double]"17.095401"
$ws.Range("I11").Value = [double]"0.0439304080539368"
$ws.Range("J11").Value = [double]"0.04393040805393679"
$ws.Range("O11").Value = [double]"0.003249135679578298"
$ws.Range("P11").Value = [double]"0.003249135679578299"
$ws.Range("Q11").Value = [double]"0.601585261701"
$ws.Range("R11").Value = [double]"5.414267355309"
$ws.Range("S11").Value = [double]"0.0001427358562264799"
$ws.Range("T11").Value = [double]"0.0001427358562264799"

# Row 12
$ws.Range("G12").Value = [double]"5.698467"
$ws.Range("H12").Value = [double]"17.095401"
$ws.Range("I12").Value = [double]"0.0439304080539368"
$ws.Range("J12").Value = [double]"0.04393040805393679"
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = [double]"0.3333333333333333"
$ws.Range("M12").Value = [double]"0.02449766666666667"
$ws.Range("N12").Value = [double]"0.073493"
$ws.Range("O12").Value = [double]"0.0007539688752111494"
$ws.Range("P12").Value = [double]"0.0007539688752111494"
$ws.Range("Q12").Value = [double]"0.139599145077"
$ws.Range("R12").Value = [double]"1.256392305693"
$ws.Range("S12").Value = [double]"3.312216034799354E-05"
$ws.Range("T12").Value = [double]"3.312216034799354E-05"

# Row 13
$ws.Range("G13").Value = [double]"5.698467"
$ws.Range("H13").Value = [double]"17.095401"
$ws.Range("I13").Value = [double]"0.0439304080539368"
$ws.Range("J13").Value = [double]"0.04393040805393679"
$ws.Range("M13").Value = [double]"32.353591"
$ws.Range("N13").Value = [double]"97.060773"
$ws.Range("O13").Value = [double]"0.9957520014958525"
$ws.Range("P13").Value = [double]"0.9957520014958525"
$ws.Range("Q13").Value = [double]"184.365870644997"
$ws.Range("R13").Value = [double]"1659.292835804973"
$ws.Range("S13").Value = [double]"0.04374379174623708"
$ws.Range("T13").Value = [double]"0.04374379174623708"

# Row 14
$ws.Range("G14").Value = [double]"0.8209666666666666"
$ws.Range("H14").Value = [double]"2.4629"
$ws.Range("I14").Value = [double]"0.006328965433220369"
$ws.Range("J14").Value = [double]"0.006328965433220369"
$ws.Range("M14").Value = [double]"0.007957000000000001"
$ws.Range("N14").Value = [double]"0.023871"
$ws.Range("O14").Value = [double]"0.0002448939493579708"
$ws.Range("P14").Value = [double]"0.0002448939493579708"
$ws.Range("Q14").Value = [double]"0.006532431766666667"
$ws.Range("R14").Value = [double]"0.05879188589999999"
$ws.Range("S14").Value = [double]"1.549925340291416E-06"
$ws.Range("T14").Value = [double]"1.549925340291416E-06"

# Row 15
$ws.Range("G15").Value = [double]"0.8209666666666666"
$ws.Range("H15").Value = [double]"2.4629"
$ws.Range("I15").Value = [double]"0.006328965433220369"
$ws.Range("J15").Value = [double]"0.006328965433220369"
$ws.Range("O15").Value = [double]"0.003249135679578298"
$ws.Range("P15").Value = [double]"0.003249135679578299"
$ws.Range("Q15").Value = [double]"0.08666917734444445"
$ws.Range("R15").Value = [double]"0.7800225961"
$ws.Range("S15").Value = [double]"2.056366740389402E-05"
$ws.Range("T15").Value = [double]"2.056366740389403E-05"

# Row 16
$ws.Range("G16").Value = [double]"0.8209666666666666"
$ws.Range("H16").Value = [double]"2.4629"
$ws.Range("I16").Value = [double]"0.006328965433220369"
$ws.Range("J16").Value = [double]"0.006328965433220369"
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = [double]"0.3333333333333333"
$ws.Range("M16").Value = [double]"0.02449766666666667"
$ws.Range("N16").Value = [double]"0.073493"
$ws.Range("O16").Value = [double]"0.0007539688752111494"
$ws.Range("P16").Value = [double]"0.0007539688752111494"
$ws.Range("Q16").Value = [double]"0.02011176774444444"
$ws.Range("R16").Value = [double]"0.1810059097"
$ws.Range("S16").Value = [double]"4.771842948935407E-06"
$ws.Range("T16").Value = [double]"4.771842948935407E-06"

# Row 17
$ws.Range("G17").Value = [double]"0.8209666666666666"
$ws.Range("H17").Value = [double]"2.4629"
$ws.Range("I17").Value = [double]"0.006328965433220369"
$ws.Range("J17").Value = [double]"0.006328965433220369"
$ws.Range("M17").Value = [double]"32.353591"
$ws.Range("N17").Value = [double]"97.060773"
$ws.Range("O17").Value = [double]"0.9957520014958525"
$ws.Range("P17").Value = [double]"0.9957520014958525"
$ws.Range("Q17").Value = [double]"26.56121975796666"
$ws.Range("R17").Value = [double]"239.0509778217"
$ws.Range("S17").Value = [double]"0.006302079997527248"
$ws.Range("T17").Value = [double]"0.006302079997527248"
